# "Added witch boss battle" - rebalance the four "levels" sheets and
# update the active-sheet/selection state left behind by the edit.

$wb = $excel.ActiveWorkbook

# --- tough_levels ----------------------------------------------------
$wsTough = $wb.Worksheets.Item("tough_levels")
$wsTough.Range("C3").Value = 100
$wsTough.Range("D3").Value = 25
$wsTough.Range("D9").Select()

# --- intelligent_levels ------------------------------------------------
$wsIntelligent = $wb.Worksheets.Item("intelligent_levels")
$wsIntelligent.Range("C3").Value = 100
$wsIntelligent.Range("D3").Value = 10
$wsIntelligent.Range("G3").Value = 20
$wsIntelligent.Range("G9").Select()

# --- wise_levels ---------------------------------------------------
$wsWise = $wb.Worksheets.Item("wise_levels")
$wsWise.Range("C3").Value = 100
$wsWise.Range("D3").Value = 10
$wsWise.Range("E3").Value = 20
$wsWise.Range("E9").Select()

# --- intuitive_levels (ends up the active tab) --------------------------
$wsIntuitive = $wb.Worksheets.Item("intuitive_levels")
$wsIntuitive.Range("C3").Value = 100
$wsIntuitive.Range("D3").Value = 10
$wsIntuitive.Range("F3").Value = 20
$wsIntuitive.Range("E27").Select()
